$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (column G) values, keyed by row number, per regenerated save_data
$kValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 1
    17 = 2
    18 = 1
    19 = 0
    21 = 2
    22 = 1
    23 = 1
    24 = 2
    25 = 2
    26 = 0
    27 = 2
    28 = 3
    29 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
